# Update status for a81f6c4a-9f7b-4245-8ab5-b0948593b50b.md from
# "Ready for handoff" to "In Translation" across all report sheets,
# as part of generating the report for archive.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B7").Value = "In Translation"
$overview.Range("C7").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B7").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B7").Value = "In Translation"
